$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# QTY <= 10 changes: clear the "Pallet Label #" (column D) values for rows 21-32
$ws.Range("D21:D32").ClearContents()
